$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-22 Friday" "2024-11-23 Saturday"

Replace-Text "115×9=1035" "349×7=2443"
Replace-Text "966×7=6762" "959×8=7672"
Replace-Text "136×5=680" "949×4=3796"
Replace-Text "232×7=1624" "301×2=602"
Replace-Text "326×3=978" "676×3=2028"

Replace-Text "856×5=4280" "915×4=3660"
Replace-Text "471×3=1413" "741×7=5187"
Replace-Text "398×2=796" "618×5=3090"
Replace-Text "662×8=5296" "770×9=6930"
Replace-Text "103×6=618" "115×2=230"

Replace-Text "169×6=1014" "754×9=6786"
Replace-Text "398×8=3184" "988×7=6916"
Replace-Text "302×5=1510" "839×3=2517"
Replace-Text "800×4=3200" "613×5=3065"
Replace-Text "911×3=2733" "627×2=1254"

Replace-Text "453×7=3171" "886×4=3544"
Replace-Text "471×8=3768" "858×4=3432"
Replace-Text "665×7=4655" "876×9=7884"
Replace-Text "887×2=1774" "397×5=1985"
Replace-Text "811×5=4055" "306×3=918"

Replace-Text "382×2=764" "216×4=864"
Replace-Text "238×3=714" "433×8=3464"
Replace-Text "499×2=998" "493×9=4437"
Replace-Text "958×4=3832" "435×9=3915"
Replace-Text "722×5=3610" "732×5=3660"
